# 392-RBI-EI-DB-DL-REC-NON-RNI-CTRFD-DL-MD-TR-1-EarlyRePayment-Newcreateloan.xlsx
# "Loan RBI, Variable Instalments"
#
# The Repayment Schedule sheet gets a new blank "Late" column inserted
# (the old "Late" / column N data shifts one column right to O, and the
# trailing "Outstanding" column shifts from P to Q). The Repayment
# Schedule tab also becomes the active/selected sheet (instead of
# NewLoanInput), with its selection moved to I21.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N - shifts old N ("Late") -> O and old P
# ("Outstanding") -> Q, leaving the new N column empty.
$wsRepay.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab, and move its selection.
$wsRepay.Select()
$wsRepay.Range("I21").Select()
